$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMAT24")
$ws.Range("A1").Value = "test"
